# Apply the change described by the diff to the FindCarTest worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FindCarTest")

# Change the runmode column (C) from "Y" to "y" for the data rows (2-4)
$ws.Range("C2").Value = "y"
$ws.Range("C3").Value = "y"
$ws.Range("C4").Value = "y"

# Update the active selection on this sheet to C5 (matches new cursor position)
$ws.Activate()
$ws.Range("C5").Select()
